# Apply "team search functionality" changes: anonymize location/country/org
# for rows 2-4 on Sheet1 (G, H, I columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

foreach ($row in 2..4) {
    $ws.Range("G$row").Value = "unknown"
    $ws.Range("H$row").Value = "Unknown"
    $ws.Range("I$row").Value = "unknown"
}
